$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 79002
$ws.Range("B3").Value = 91809
$ws.Range("B4").Value = 92107
$ws.Range("B5").Value = 93096
$ws.Range("B6").Value = 79244
$ws.Range("B7").Value = 92107
$ws.Range("B8").Value = 92107
